# Apply the recorded Solver run: add the hidden solver_* defined names
# (sheet-scoped) that Excel's Solver add-in stores after an optimization,
# update the decision-variable weights (C15:G15) with the new, correlation
# -free optimum, and move the active selection to L13 as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Solver parameters (hidden, sheet-scoped defined names) ---------------
$solverNames = @(
    @('solver_adj',  'Calcs!$C$15:$G$15'),
    @('solver_cvg',  '0.0001'),
    @('solver_drv',  '1'),
    @('solver_eng',  '1'),
    @('solver_est',  '1'),
    @('solver_itr',  '2147483647'),
    @('solver_lhs1', 'Calcs!$C$6'),
    @('solver_lhs2', 'Calcs!$H$15'),
    @('solver_mip',  '2147483647'),
    @('solver_mni',  '30'),
    @('solver_mrt',  '0.075'),
    @('solver_msl',  '2'),
    @('solver_neg',  '1'),
    @('solver_nod',  '2147483647'),
    @('solver_num',  '2'),
    @('solver_nwt',  '1'),
    @('solver_opt',  'Calcs!$C$3'),
    @('solver_pre',  '0.000001'),
    @('solver_rbv',  '1'),
    @('solver_rel1', '2'),
    @('solver_rel2', '2'),
    @('solver_rhs1', '0.05'),
    @('solver_rhs2', '1'),
    @('solver_rlx',  '2'),
    @('solver_rsd',  '0'),
    @('solver_scl',  '1'),
    @('solver_sho',  '2'),
    @('solver_ssz',  '100'),
    @('solver_tim',  '2147483647'),
    @('solver_tol',  '0.01'),
    @('solver_typ',  '2'),
    @('solver_val',  '0'),
    @('solver_ver',  '3')
)

foreach ($pair in $solverNames) {
    $defName = $ws.Names.Add($pair[0], '=' + $pair[1])
    $defName.Visible = $false
}

# --- Updated optimal weights (Solver result w/out correlation) ------------
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.41085053811994493
$ws.Range("F15").Value = -0.00000048434121259268236
$ws.Range("G15").Value = 0.58914993859401388

# --- Restore the author's last selection -----------------------------------
$ws.Range("L13").Select()
